$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trang_tính1")

# Move the assigned-person / task text from the "Day 2" column (D) to the
# "Day 4" column (H) for rows 16-20, matching the updated weekly plan
# (assignments shifted to 12/10/22).
for ($r = 16; $r -le 20; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $dCell.Value2
    $dCell.Value = $null
}

$ws.Range("I16").Select()

# Reflect the view state left behind after the edit (scrolled down a bit,
# zoomed in from 60% to 70%).
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
